# Insert a new data row above current row 38 (shifts existing rows 38-44
# down to 39-45, preserving all their values/formatting), then populate the
# newly inserted row 38 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38; this pushes rows 38:44 down to 39:45.
$ws.Rows.Item(38).Insert()

# Fill the new row 38 with the new record's data.
$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"

$fecha38 = Get-Date -Year 2023 -Month 1 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Range("D38").Value = $fecha38

$ws.Range("E38").Value = 10
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100101
$ws.Range("H38").Value = "Berries"
$ws.Range("I38").Value = 100101001
$ws.Range("J38").Value = "Arándano (blue)"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 2000
$ws.Range("O38").Value = 2200
$ws.Range("P38").Value = 2100
$ws.Range("Q38").Value = '$/bandeja 2 kilos'
$ws.Range("R38").Value = "Provincia de Colchagua"
$ws.Range("S38").Value = 1050
$ws.Range("T38").Value = 2
